$d = $word.ActiveDocument

$d.Content.Find.Execute("2025-06-09 Monday", $true, $false, $false, $false, $false, $true, 1, $false, "2025-06-10 Tuesday", 2) | Out-Null
$d.Content.Find.Execute("814×9=", $true, $false, $false, $false, $false, $true, 1, $false, "880×3=", 2) | Out-Null
$d.Content.Find.Execute("259×6=", $true, $false, $false, $false, $false, $true, 1, $false, "831×8=", 2) | Out-Null
$d.Content.Find.Execute("407×8=", $true, $false, $false, $false, $false, $true, 1, $false, "687×2=", 2) | Out-Null
$d.Content.Find.Execute("541×5=", $true, $false, $false, $false, $false, $true, 1, $false, "621×6=", 2) | Out-Null
$d.Content.Find.Execute("974×9=", $true, $false, $false, $false, $false, $true, 1, $false, "474×9=", 2) | Out-Null
$d.Content.Find.Execute("390×7=", $true, $false, $false, $false, $false, $true, 1, $false, "834×4=", 2) | Out-Null
$d.Content.Find.Execute("510×6=", $true, $false, $false, $false, $false, $true, 1, $false, "707×2=", 2) | Out-Null
$d.Content.Find.Execute("736×2=", $true, $false, $false, $false, $false, $true, 1, $false, "414×8=", 2) | Out-Null
$d.Content.Find.Execute("441×7=", $true, $false, $false, $false, $false, $true, 1, $false, "861×6=", 2) | Out-Null
$d.Content.Find.Execute("205×8=", $true, $false, $false, $false, $false, $true, 1, $false, "950×5=", 2) | Out-Null
$d.Content.Find.Execute("919×3=", $true, $false, $false, $false, $false, $true, 1, $false, "408×5=", 2) | Out-Null
$d.Content.Find.Execute("843×6=", $true, $false, $false, $false, $false, $true, 1, $false, "628×3=", 2) | Out-Null
$d.Content.Find.Execute("334×9=", $true, $false, $false, $false, $false, $true, 1, $false, "801×2=", 2) | Out-Null
$d.Content.Find.Execute("946×7=", $true, $false, $false, $false, $false, $true, 1, $false, "755×7=", 2) | Out-Null
$d.Content.Find.Execute("763×6=", $true, $false, $false, $false, $false, $true, 1, $false, "575×3=", 2) | Out-Null
$d.Content.Find.Execute("864×5=", $true, $false, $false, $false, $false, $true, 1, $false, "797×5=", 2) | Out-Null
$d.Content.Find.Execute("674×5=", $true, $false, $false, $false, $false, $true, 1, $false, "533×6=", 2) | Out-Null
$d.Content.Find.Execute("388×3=", $true, $false, $false, $false, $false, $true, 1, $false, "750×8=", 2) | Out-Null
$d.Content.Find.Execute("326×2=", $true, $false, $false, $false, $false, $true, 1, $false, "394×3=", 2) | Out-Null
$d.Content.Find.Execute("498×5=", $true, $false, $false, $false, $false, $true, 1, $false, "851×8=", 2) | Out-Null
$d.Content.Find.Execute("125×2=", $true, $false, $false, $false, $false, $true, 1, $false, "317×4=", 2) | Out-Null
$d.Content.Find.Execute("201×2=", $true, $false, $false, $false, $false, $true, 1, $false, "983×9=", 2) | Out-Null
$d.Content.Find.Execute("231×6=", $true, $false, $false, $false, $false, $true, 1, $false, "218×6=", 2) | Out-Null
$d.Content.Find.Execute("252×7=", $true, $false, $false, $false, $false, $true, 1, $false, "261×3=", 2) | Out-Null
$d.Content.Find.Execute("856×5=", $true, $false, $false, $false, $false, $true, 1, $false, "300×4=", 2) | Out-Null
